$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B5").Value = "'-1000000000000000000000000000 .. 2000"
